$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16
$ws.Cells.Item(16,1).Value = 'Informatica Cloud Data Governance Catalog Specialist'
$ws.Cells.Item(16,2).Value = 'https://www.dice.com/job-detail/6f16b244-a92c-425a-9007-b06701ccabcf?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang'
$ws.Cells.Item(16,3).Value = 'Torrance, California'
$ws.Cells.Item(16,4).Value = 'Contract'
$ws.Cells.Item(16,5).Value = 'USD 53.00 - 55.00 per hour'
$ws.Cells.Item(16,6).Value = 'Robert Half'
$ws.Cells.Item(16,7).Value = 'Description We are looking for an experienced Informatica Cloud Data Governance Catalog Specialist to join our team in Southern California. This position involves working on-site four days per week and offers a long-term contract opportunity. The ideal candidate will have a strong background in data governance, analytics, and business intelligence tools, coupled with a proactive approach to problem-solving and collaboration. Responsibilities: Create catalog quality reports to monitor and enhance data governance metrics across domains and sub-domains. Develop and showcase data governance dashboards tailored to different user roles, including Data Owners, Stewards, Engineers, and Privacy Officers. Collaborate with business and IT teams, including data stewards, catalog architects, and platform owners, to implement governance solutions. Execute profiling, sampling, and scanner setups using Informatica tools to ensure data quality. Apply expertise in metadata management, data modeling, and large-scale data analysis to support governance initiatives. Design and implement both traditional relational and modern big-data architectures based on organizational requirements. Utilize business intelligence tools such as Power BI and Tableau to create actionable insights and reports. Define compliance procedures and produce audit reports to meet regulatory requirements. Establish and support governance councils and operational frameworks using data catalog tools. Facilitate metadata ingestion and ensure adherence to data security and quality standards. Requirements Proficiency with tools like Informatica Cloud Data Governance Catalog and Cloud Data Quality. Hands-on experience in data modeling, metadata management, and large-scale data analysis. Familiarity with Collibra, Alation, and Glue Data Catalog. Strong understanding of entity-relationship modeling and data security practices. Expertise in business intelligence technologies such as Power BI and Tableau. Exceptional communication and presentation skills to effectively convey technical concepts. Analytical mindset with proven problem-solving abilities. Ability to work collaboratively as part of a team and build strong relationships with stakeholders. Technology Doesn''t Change the World, People Do. Robert Half is the world''s first and largest specialized talent solutions firm that connects highly qualified job seekers to opportunities at great companies. We offer contract, temporary and permanent placement solutions for finance and accounting, technology, marketing and creative, legal, and administrative and customer support roles. Robert Half works to put you in the best position to succeed. We provide access to top jobs, competitive compensation and benefits, and free online training. Stay on top of every opportunity - whenever you choose - even on the go. Download the Robert Half app and get 1-tap apply, notifications of AI-matched jobs, and much more. All applicants applying for U.S. job openings must be legally authorized to work in the United States. Benefits are available to contract/temporary professionals, including medical, vision, dental, and life and disability insurance. Hired contract/temporary professionals are also eligible to enroll in our company 401(k) plan. Visit roberthalf.gobenefits.net for more information. 2025 Robert Half. An Equal Opportunity Employer. M/F/Disability/Veterans. By clicking "Apply Now," you''re agreeing to Robert Half''s Terms of Use and Privacy Notice.'
$ws.Cells.Item(16,8).Value = '24.4%%'
$ws.Cells.Item(16,9).Value = 'Contract Independent | Contract W2 | On-site | USD $53.00 - 55.00 per hour'

# Row 17
$ws.Cells.Item(17,1).Value = 'Java Integration Architect - W2 Contract - Ongoing'
$ws.Cells.Item(17,2).Value = 'https://www.dice.com/job-detail/23ae0c27-91c5-4fc5-aa7c-90aaa609844f?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=2&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang'
$ws.Cells.Item(17,3).Value = 'Hybrid in Oakland, California'
$ws.Cells.Item(17,4).Value = 'Contract'
$ws.Cells.Item(17,5).Value = '$75 - $85'
$ws.Cells.Item(17,6).Value = 'Triune Infomatics Inc'
$ws.Cells.Item(17,7).Value = 'Role: Java Integration Architect Location: Oakland, CA (2 days onsite mandatory) Duration: 3 months + extensions >> ongoing Overview: Lead a 3 months'' discovery, design, implementation support and remediation exercise so existing on-premises Java applications (Linux/JBoss) can interoperate with workloads migrating to cloud. Inspect code and runtime configuration to identify every change needed to repoint integrations and produce concise compatibility patterns to resolve common issues through testing and go-live. A major stream is re-platforming DB2 integrations to Microsoft SQL Server, covering both configuration and SQL/dialect changes. Project Insights: DB2 to SQL Server migration will involve confirmed complexities such as stored procedures; additional challenges like schema differences, SQL dialect variations, and data type mismatches are expected during delivery. Architect will collaborate with the existing DBA to analyze current models and resolve issues related to embedded SQL logic, pagination, identity/sequence handling, and stored procedure calls. Migration from on-prem DB2 to cloud SQL Server will require assessment and adaptation of TLS configurations, identity management practices, and network/proxy constraints — exact impacts to be identified as part of solution discovery. Key Responsibilities: Discovery & code/config review: Analyze Java apps on Linux/JBoss to locate all integration touchpoints (HTTP, messaging, databases, files, schedulers, security, network paths) with exact file/line, JNDI and environment references. Application Remediation Report (deliverable): Capture the current → target change set, change action, and validation steps for each Java application. Examples include: Connectivity : New endpoints/FQDNs/ports for migrated APIs/gateways/queues/databases; DNS/TTL; routing/peering/private link; proxy/NO_PROXY where required; DR/failover targets for those endpoints. Security : TLS specifics for new endpoints (server names/SANs, cert chains, trust/keystores, ciphers, hostname verification/SNI); identities/secrets used to call the migrated services (users/passwords, API keys/tokens), where stored, and rotation owner. Interfaces : HTTP APIs now fronted in cloud: base URLs, versions/paths, required headers/scopes, contract deltas and version pinning. Gateways referenced by Java (only if used): e.g., CTG/cloud API gateway host/port/alias/timeouts/pools/charset as called by the app. Messaging (JMS/MQ) only if that queueing endpoint is migrating: queue manager/channel/connection names, ports, TLS, credentials. Database : New driver and JDBC URL (host/port/db, encryption, failover/AG listener/multiSubnetFailover), credentials, and required SQL/dialect/schema changes (types, functions, pagination, isolation/locking/WITH UR equivalents, stored procedure call syntax, error/return code mapping, identifier/collation rules). Include the network path to SQL Server (private link/peering/proxy) and validation of connectivity/failover. Runtime & resilience : Connect/read timeouts, retries/back-off (idempotency where relevant), circuit breakers/health checks, transaction recovery impacts; charset/locale where gateway/database behavior changes. Application Compatibility Patterns (deliverable) : Produce short, reusable playbooks within Java scope for remediation and topology shifts (e.g., JDBC→SQL Server, TLS trust/hostname checks, HTTP client timeouts/retries, MQ/JMS client config, proxy/NO_PROXY, DNS failover, charset; and Java→CTG→CICS to Java→Cloud API Gateway→service; Java→VPN→DB2 to Java→SQL Server via private link), each with minimal code/config, verification and rollback. Stakeholder engagement: Work with client architects (application, integration, security, network) to shape options, agree patterns and support UAT/cutover. Governance & traceability: Maintain per-app change logs and risks; surface dependencies and non-functional impacts (latency, throughput, DR). Documentation: Design and document the implementation steps to guide the development team through execution. Update and revise the documentation to align with the final integration setup to train and hand over the maintainability of the integrations to the client support team. Skills & Experience (Essential): 7+ years as an application architect or senior Java engineer, able to read Java and pinpoint integration/configuration changes in Linux/JBoss environments. Demonstrable delivery of discovery/remediation for on-prem → cloud migrations, re-pointing integrations. Strong SQL with proven Microsoft SQL Server experience, including assessing/adapting SQL and schema when moving from DB2. Solid understanding in networking/DNS, routing/proxies and firewall/security-group models for app connectivity. Practical TLS/PKI knowledge for Java/JBoss (keystores/truststores, hostname verification, etc.). Clear technical writing and stakeholder management; comfortable operating independently while supporting client teams. Desirable: Familiarity with IBM MQ/JMS client configuration from Java. Exposure to feature-flag/remote-config approaches. Experience with common Java web stacks (Servlet/JAX-RS; Spring familiarity is a plus but not required). Experience with DR testing and performance tuning under changing latency profiles.'
$ws.Cells.Item(17,8).Value = '35.79%%'
$ws.Cells.Item(17,9).Value = 'Contract W2 | No Travel Required | Hybrid | $75 - $85/hr'
